# Regenerate the "K" column (column G) values for save_data sheet.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the visible effect on this worksheet is that the
# per-row values in column G (header "K") are recalculated/rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K"), taken from the
# regenerated save_data output.
$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 3
    15 = 2
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 2
    29 = 0
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
